$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = "ProductIdTest1"
$ws.Range("Q2").Value = "ProductIdTest2"

$ws.Range("P3").Value = "C2EE3694-959A-4A87-BC8C-4003F6576352"
$ws.Range("P4").Value = "C2EE3694-959A-4A87-BC8C-4003F6576353"
$ws.Range("P5").Value = "C2EE3694-959A-4A87-BC8C-4003F6576354"
$ws.Range("P6").Value = "C2EE3694-959A-4A87-BC8C-4003F6576355"
$ws.Range("P7").Value = "C2EE3694-959A-4A87-BC8C-4003F6576356"
$ws.Range("P8").Value = "C2EE3694-959A-4A87-BC8C-4003F6576357"

$ws.Range("Q3").Value = "C2EE3694-959A-4A87-BC8C-4003F6576357"
$ws.Range("Q4").Value = "C2EE3694-959A-4A87-BC8C-4003F6576358"
$ws.Range("Q5").Value = "C2EE3694-959A-4A87-BC8C-4003F6576359"

$ws.Columns.Item(16).ColumnWidth = 39.285714285714285
$ws.Columns.Item(17).ColumnWidth = 37.285714285714285

[void]$ws.Range("Q11").Select()
